$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting all existing data right by one.
$ws.Range("A1").EntireColumn.Insert()

# Add header for the new leaderID column.
$ws.Range("A1").Value = "leaderID"

# Fill leaderID values 0..19 for the 20 data rows (now rows 2-21).
for ($i = 0; $i -le 19; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $i
}

# The filter database range needs to shift right along with the data
# (column insert does not retarget named ranges automatically).
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=Sheet1!`$B`$1:`$F`$130"
    }
}

# Match the new active selection cell from the diff.
$ws.Range("B24").Select() | Out-Null
